$p = $ppt.ActivePresentation

# Slide 1
$s = $p.Slides.Item(1)

# Title (shape 1): replace text without leaving stray rPr/paragraph artifacts
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "X`r微軟企業文化 - 第一頁"
$titleTr.Paragraphs(1).Delete()

# Content placeholder (shape 2): replace the 5 bullet paragraphs
$bodyTr = $s.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "X`r• 微軟致力於透過科技引領全球進步。`r• 我們的核心價值強調包容及多樣性。`r• 社會責任是企業文化的基石。`r• 與客戶和合作伙伴密切協作是我們的理念。`r• 微軟支持創新，並鼓勵大膽嘗試。"
$bodyTr.Paragraphs(1).Delete()

# Slide 2
$s = $p.Slides.Item(2)

# Title (shape 1): replace text without leaving stray rPr/paragraph artifacts
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "X`r微軟企業文化 - 第二頁"
$titleTr.Paragraphs(1).Delete()

# Content placeholder (shape 2): replace the 5 bullet paragraphs
$bodyTr = $s.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "X`r• 微軟的工作環境注重開放性與透明度。`r• 員工享有彈性工作安排並支持遠端作業。`r• 對技術的熱情是我們不斷前進的驅動力。`r• 領導支持員工的成長與學習。`r• 微軟追求卓越，並實現可持續性目標。"
$bodyTr.Paragraphs(1).Delete()

# Slide 3
$s = $p.Slides.Item(3)

# Title (shape 1): replace text without leaving stray rPr/paragraph artifacts
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "X`r微軟企業文化 - 第三頁"
$titleTr.Paragraphs(1).Delete()

# Content placeholder (shape 2): replace the 5 bullet paragraphs
$bodyTr = $s.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "X`r• 微軟鼓勵員工參與創新項目。`r• 我們的企業文化以信任為基礎。`r• 顧客的成功是我們的成功。`r• 微軟重視社群及環境影響。`r• 推動技術教育以建立更好的未來。"
$bodyTr.Paragraphs(1).Delete()

# Slide 4
$s = $p.Slides.Item(4)

# Title (shape 1): replace text without leaving stray rPr/paragraph artifacts
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "X`r微軟企業文化 - 第四頁"
$titleTr.Paragraphs(1).Delete()

# Content placeholder (shape 2): replace the 5 bullet paragraphs
$bodyTr = $s.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "X`r• 微軟支持包容性，建立多元化的技術行業。`r• 我們的產品是以用戶需求為中心設計的。`r• 微軟關注環境保護並實現碳中和目標。`r• 員工的健康與幸福是文化的一部分。`r• 微軟的使命是賦能個人及機構創造更多成果。"
$bodyTr.Paragraphs(1).Delete()

# Slide 5
$s = $p.Slides.Item(5)

# Title (shape 1): replace text without leaving stray rPr/paragraph artifacts
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "X`r微軟企業文化 - 第五頁"
$titleTr.Paragraphs(1).Delete()

# Content placeholder (shape 2): replace the 5 bullet paragraphs
$bodyTr = $s.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "X`r• 微軟的品牌理念體現信任與可靠性。`r• 我們鼓勵員工擁有創造性思維並挑戰自我。`r• 微軟與非營利組織合作以提高社會影響力。`r• 投資技術以幫助全球應對重要挑戰。`r• 微軟致力於提供世界一流的解決方案。"
$bodyTr.Paragraphs(1).Delete()
